# FTT-S: Classification Titles update 2_04062024
# Added SMTI, SSTI, XTPI

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the three new classification-title sheets, right after C5TI, in
#    SMTI -> SSTI -> XTPI order (matches sheetId 37/38/39, rId29/30/31).
# ---------------------------------------------------------------------------
$c5ti = $wb.Worksheets.Item("C5TI")

$smti = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $c5ti)
$smti.Name = "SMTI"

$ssti = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $smti)
$ssti.Name = "SSTI"

$xtpi = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ssti)
$xtpi.Name = "XTPI"

# ---------------------------------------------------------------------------
# 2. Populate SMTI (Steel Material Titles)
# ---------------------------------------------------------------------------
$smtiNames = @(
    "Crude Steel",
    "Primary Iron Production (PIP)",
    "Primary Iron Import (PII)",
    "Scrap",
    "Coke",
    "Biocoke",
    "Sinter",
    "Sinter (BB)",
    "Pellets",
    "Pellets (BB)",
    "Oxygen",
    "Hard Coal",
    "Other Coal",
    "Natural Gas",
    "Electricity",
    "Iron ore",
    "Limestone",
    "Hydrogen",
    "Biocharcoal",
    "Biogas",
    "IC",
    "O&M",
    "EF",
    "Employment"
)

$smti.Range("A1").Value = "Full name"
$smti.Range("B1").Value = "Short name"
for ($i = 0; $i -lt $smtiNames.Count; $i++) {
    $row = $i + 2
    $smti.Cells.Item($row, 1).Value = $smtiNames[$i]
    $smti.Cells.Item($row, 2).Value = $i + 1
}
$smti.Columns.Item(1).ColumnWidth = 11.5

# ---------------------------------------------------------------------------
# 3. Populate SSTI (Steel Sector Titles)
# ---------------------------------------------------------------------------
$sstiNames = @(
    "Coke plant",
    "Biocoke plant",
    "Sinter plant",
    "Biosinter plant",
    "Pellet plant",
    "Biopellet plant",
    "Oxygen plant",
    "Conv. BF",
    "Conv. BF (BB)",
    "BF TGR (CCS)",
    "BF TGR (CCS, BB)",
    "DR-gas",
    "DR-gas (BB)",
    "DR-coal",
    "DR-coal (BB)",
    "SR",
    "SR (BB)",
    "SR+ ",
    "SR+ (BB)",
    "HFS",
    "OHF",
    "BOF",
    "BOF (BB)",
    "EAF",
    "EAF (BB)",
    "MOE",
    "Alkaline electrolysis",
    "HFS",
    "Final stage"
)

$ssti.Range("A1").Value = "Full name"
$ssti.Range("B1").Value = "Short name"
for ($i = 0; $i -lt $sstiNames.Count; $i++) {
    $row = $i + 2
    $ssti.Cells.Item($row, 1).Value = $sstiNames[$i]
    $ssti.Cells.Item($row, 2).Value = $i + 1
}
$ssti.Columns.Item(1).ColumnWidth = 18.3

# ---------------------------------------------------------------------------
# 4. Populate XTPI (eXTra Production-stage Indices)
# ---------------------------------------------------------------------------
$xtpiNames = @(
    "Transport ",
    "Machinery",
    "Construction",
    "Products"
)

$xtpi.Range("A1").Value = "Full name"
$xtpi.Range("B1").Value = "Short name"
for ($i = 0; $i -lt $xtpiNames.Count; $i++) {
    $row = $i + 2
    $xtpi.Cells.Item($row, 1).Value = $xtpiNames[$i]
    $xtpi.Cells.Item($row, 2).Value = $i + 1
}
$xtpi.Columns.Item(1).ColumnWidth = 26.5

# ---------------------------------------------------------------------------
# 5. Make SSTI the active/selected sheet (matches activeTab=29 / tabSelected
#    moving off STTI onto SSTI in the target workbook).
# ---------------------------------------------------------------------------
$ssti.Select()
